$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(45954, 5594, 3990, 3673, 233, 46, 35, 3, 0),
    @(45957, 5592, 4178, 3795, 291, 60, 25, 4, 3)
)

$startRow = 14
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $vals = $data[$i]
    for ($c = 0; $c -lt $vals.Count; $c++) {
        $ws.Cells.Item($row, $c + 1).Value = $vals[$c]
    }
}

$ws.Range("G18").Select()
